$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data rows were reordered (one JSON record's worth of
# time-bucket analysis data rotated in): what was row 3 becomes row 2,
# what was row 4 becomes row 3, and what was row 2 becomes row 4.
# Hyperlink click-targets (the rId relationships) stay pinned to their
# original row position; only the displayed text (title/timestamp/
# historical distance/uri) moves with the logical row.

$ws.Range("A2").Value = "Ebola outbreak: Red Cross workers attacked while burying dead bodies in Guinea"
$ws.Range("B2").Value = "2014-09-24T20:27:42UTC"
$ws.Range("C2").Value = 23
$ws.Range("D2").Value = "day_2_to_30"
$ws.Range("E2").Value = "https://www.independent.co.uk/news/world/africa/ebola-outbreak-red-cross-workers-attacked-while-burying-dead-bodies-in-guinea-9754140.html"

$ws.Range("A3").Value = "Journalists, health team killed while conducting Ebola awareness-raising campaign"
$ws.Range("B3").Value = "2014-09-19T00:00:00UTC"
$ws.Range("C3").Value = 18
$ws.Range("D3").Value = "day_2_to_30"
$ws.Range("E3").Value = "https://www.ifex.org/guinea/2014/09/19/journalists_killed/"

$ws.Range("A4").Value = "Guinea arrests 27 over Ebola health team murders"
$ws.Range("B4").Value = "2014-09-24T00:00:00UTC"
$ws.Range("C4").Value = 23
$ws.Range("D4").Value = "day_2_to_30"
$ws.Range("E4").Value = "http://reliefweb.int/report/guinea/guinea-arrests-27-over-ebola-health-team-murders"
